# Update cryptos list with latest price and 1h volume percentage values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "65.583.07";  E = "  -1.68%  " },
    @{ Row = 3;  D = "3.398.02";   E = "  -2.34%  " },
    @{ Row = 4;  D = "1.00";       E = "  -0.13%  " },
    @{ Row = 5;  D = "596.42";     E = "  -1.24%  " },
    @{ Row = 6;  D = "142.07";     E = "  -4.16%  " },
    @{ Row = 7;  D = $null;        E = "  -0.16%  " },
    @{ Row = 8;  D = "3.396.83";   E = "  -2.31%  " },
    @{ Row = 9;  D = $null;        E = "  -2.77%  " },
    @{ Row = 10; D = "7.89";       E = "  +4.13%  " },
    @{ Row = 11; D = $null;        E = "  -6.30%  " },
    @{ Row = 12; D = "0.404";      E = "  -4.84%  " },
    @{ Row = 13; D = "3.972.97";   E = "  -2.35%  " },
    @{ Row = 14; D = "0.0000199";  E = "  -7.04%  " },
    @{ Row = 15; D = "29.43";      E = "  -6.86%  " },
    @{ Row = 16; D = $null;        E = "  -0.54%  " },
    @{ Row = 17; D = "65.583.40";  E = "  -1.87%  " },
    @{ Row = 18; D = "3.393.91";   E = "  -2.80%  " },
    @{ Row = 19; D = "10.32";      E = "  +2.14%  " },
    @{ Row = 20; D = "6.09";       E = "  -5.85%  " },
    @{ Row = 21; D = "14.54";      E = "  -5.54%  " },
    @{ Row = 22; D = "413.61";     E = "  -5.89%  " },
    @{ Row = 23; D = "0.578";      E = "  -5.35%  " },
    @{ Row = 24; D = "76.97";      E = "  -3.33%  " },
    @{ Row = 25; D = $null;        E = "  +0.11%  " },
    @{ Row = 26; D = "3.532.35";   E = "  -2.41%  " },
    @{ Row = 27; D = "0.0000108";  E = "  -9.83%  " },
    @{ Row = 28; D = "9.19";       E = "  -5.86%  " },
    @{ Row = 29; D = $null;        E = "  -7.43%  " },
    @{ Row = 30; D = $null;        E = "  -3.23%  " },
    @{ Row = 31; D = $null;        E = "  -0.52%  " },
    @{ Row = 32; D = "0.159";      E = "  -5.54%  " },
    @{ Row = 33; D = $null;        E = "  -8.41%  " },
    @{ Row = 34; D = "24.36";      E = "  -4.17%  " },
    @{ Row = 35; D = "3.393.62";   E = "  -2.19%  " },
    @{ Row = 37; D = "5.50";       E = "  -8.97%  " },
    @{ Row = 38; D = $null;        E = "  -7.42%  " },
    @{ Row = 39; D = "7.48";       E = "  -5.66%  " },
    @{ Row = 40; D = "0.998";      E = "  -0.21%  " },
    @{ Row = 41; D = "168.63";     E = "  -4.16%  " },
    @{ Row = 42; D = "0.0849";     E = "  -4.47%  " },
    @{ Row = 43; D = "0.869";      E = "  -2.47%  " },
    @{ Row = 44; D = "5.01";       E = "  -7.49%  " },
    @{ Row = 45; D = "1.92";       E = "  -10.28%  " },
    @{ Row = 46; D = "45.33";      E = "  -2.03%  " },
    @{ Row = 47; D = "26.34";      E = "  -8.98%  " },
    @{ Row = 48; D = $null;        E = "  -5.08%  " },
    @{ Row = 49; D = "7.02";       E = "  -6.11%  " },
    @{ Row = 50; D = "2.24";       E = "  -7.80%  " },
    @{ Row = 51; D = "0.914";      E = "  -7.01%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($u.D -ne $null) {
        $dCell = $ws.Cells.Item($r, 4)
        # Force the price column to remain plain text even when the new
        # value looks numeric (e.g. "1.00", "596.42"): mark the cell as
        # Text before assigning, then drop the NumberFormat residue so the
        # cell's style returns to its original (unstyled) state.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.ClearFormats()
    }

    # The volume column is always textual ("  -1.68%  " etc. never parses
    # as a number because of the leading/trailing spaces and percent
    # sign), so a direct assignment is safe.
    $ws.Cells.Item($r, 5).Value = $u.E
}
